$d = $word.ActiveDocument

# --- 1. Normalise the whole document: round-tripping the body through
#        WordOpenXML / InsertXML merges adjacent same-formatted runs and
#        drops the w:proofErr spell/grammar-check markers that used to
#        split them - exactly what the target diff shows happening in
#        the "by", "Prodi", "Quote's", and "Se-Kata..." paragraphs.
$bodyXml = $d.Content.WordOpenXML
$d.Content.InsertXML($bodyXml)

# --- 1b. Drop the stray "/" between "Hasiholan" and "Omega Purba" on the
#         "by : ..." byline.
$find = $d.Content.Find
$find.Execute("by : Samuel Hasiholan /Omega Purba, S. Tr. T.", $true, $false, $false, $false, $false, $true, 1, $false, "by : Samuel Hasiholan Omega Purba, S. Tr. T.", 2) | Out-Null

# --- 2. Change the two math exponents from "x" to "2".
for ($i = 1; $i -le $d.OMaths.Count; $i++) {
    $om = $d.OMaths.Item($i)
    $rng = $om.Range
    $xml = $rng.WordOpenXML
    $newXml = $xml.Replace("<m:t>x</m:t>", "<m:t>2</m:t>")
    if ($newXml -ne $xml) {
        $rng.InsertXML($newXml)
    }
}

# --- 3. Split the founder paragraph: insert a new paragraph right after it
#        that contains "[Proverb 11 : 24]".
$founderPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Samuel Hasiholan Omega, S. Tr. T. (Founder : BeruangLaut.ID)") {
        $founderPara = $d.Paragraphs.Item($i)
        break
    }
}

$endOfPara = $founderPara.Range.End
$insertPoint = $d.Range($endOfPara, $endOfPara)
$insertPoint.InsertParagraphAfter()

# find the (now empty) paragraph that was just created and put the text in it
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Samuel Hasiholan Omega, S. Tr. T. (Founder : BeruangLaut.ID)") {
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = "[Proverb 11 : 24]"
        break
    }
}

Write-Host "Edit complete"
